$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (row 2 through row 31) from 2023-10-25 (45224) to 2023-11-03 (45233)
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 3).Value = 45233
}
